$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update party-name strings (rows 1, 2 and 4) ---
# These assignments replace the old "Scottish ..." party names with the
# shortened official names; Excel will garbage-collect the now-unused
# shared strings and append the new ones at the end of the table.
$ws.Range("C1").Value = "Conservative and Unionist Party"
$ws.Range("C2").Value = "Liberal Democrats"
$ws.Range("C4").Value = "Labour Party"

# --- Add "Count 3" formulas in column G (rows 9-14) ---
$ws.Range("G9").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D1&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B1&""",""Firstname"":"""&$A1&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C1&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""3"",""Transfers"":"""&G1&""",""id"":"&ROW()+6&",""Total_Votes"":"""&$H1&"""},"'
$ws.Range("G10:G14").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D2&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B2&""",""Firstname"":"""&$A2&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C2&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""3"",""Transfers"":"""&G2&""",""id"":"&ROW()+6&",""Total_Votes"":"""&$H2&"""},"'

# --- Add "Count 5" formulas in column K (rows 9-14) ---
$ws.Range("K9").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D1&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B1&""",""Firstname"":"""&$A1&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C1&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""5"",""Transfers"":"""&K1&""",""id"":"&ROW()+12&",""Total_Votes"":"""&L1&"""},"'
$ws.Range("K10:K14").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D2&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B2&""",""Firstname"":"""&$A2&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C2&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""5"",""Transfers"":"""&K2&""",""id"":"&ROW()+6&",""Total_Votes"":"""&L2&"""},"'

# --- Update the saved selection to match the newly added range ---
$ws.Range("K9:K14").Select()
